$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the newly published day's spot prices (automatic update)
$ws.Range("A2").Value = 45986

$ws.Range("B2").Value = 42.92
$ws.Range("C2").Value = 39.75
$ws.Range("D2").Value = 29.84
$ws.Range("E2").Value = 21.51
$ws.Range("F2").Value = 16.71
$ws.Range("G2").Value = 31.83
$ws.Range("H2").Value = 59.78
$ws.Range("I2").Value = 77.68000000000001
$ws.Range("J2").Value = 82.12
$ws.Range("K2").Value = 58.04
$ws.Range("L2").Value = 13.46
$ws.Range("M2").Value = 3.24
$ws.Range("N2").Value = 0.96
$ws.Range("O2").Value = 0.66
$ws.Range("P2").Value = 0.66
$ws.Range("Q2").Value = 2.74
$ws.Range("R2").Value = 27.56
$ws.Range("S2").Value = 89.13
$ws.Range("T2").Value = 99.83
$ws.Range("U2").Value = 102.11
$ws.Range("V2").Value = 96.28
$ws.Range("W2").Value = 90.31
$ws.Range("X2").Value = 88.59999999999999
$ws.Range("Y2").Value = 80.38
$ws.Range("Z2").Value = 48.17

# AA2 (Slot_4h_max) label is unchanged: "20h-24h"
$ws.Range("AB2").Value = 88.89

# Slot_2h_frist / Slot_2h_second swapped labels and updated prices
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 100.97
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 93.3

$ws.Range("AG2").Value = "0h-16h"
